# Refresh the crypto price/volume snapshot (columns D and E, rows 2-51) with the
# latest values pulled by the scheduled GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Text = '42.768.08'; Numeric = $false },
    @{ Cell = 'E2'; Text = '  -1.01%  '; Numeric = $false },
    @{ Cell = 'D3'; Text = '2.370.37'; Numeric = $false },
    @{ Cell = 'E3'; Text = '  +1.06%  '; Numeric = $false },
    @{ Cell = 'E4'; Text = '  -0.19%  '; Numeric = $false },
    @{ Cell = 'D5'; Text = '331.98'; Numeric = $true },
    @{ Cell = 'E5'; Text = '  +5.82%  '; Numeric = $false },
    @{ Cell = 'D6'; Text = '101.20'; Numeric = $true },
    @{ Cell = 'E6'; Text = '  -7.71%  '; Numeric = $false },
    @{ Cell = 'E7'; Text = '  -0.38%  '; Numeric = $false },
    @{ Cell = 'E8'; Text = '  +0.02%  '; Numeric = $false },
    @{ Cell = 'D9'; Text = '0.630'; Numeric = $true },
    @{ Cell = 'E9'; Text = '  -0.36%  '; Numeric = $false },
    @{ Cell = 'D10'; Text = '40.10'; Numeric = $true },
    @{ Cell = 'E10'; Text = '  -6.48%  '; Numeric = $false },
    @{ Cell = 'D11'; Text = '0.0925'; Numeric = $true },
    @{ Cell = 'E11'; Text = '  -1.41%  '; Numeric = $false },
    @{ Cell = 'D12'; Text = '8.50'; Numeric = $true },
    @{ Cell = 'E12'; Text = '  -4.12%  '; Numeric = $false },
    @{ Cell = 'E13'; Text = '  -3.39%  '; Numeric = $false },
    @{ Cell = 'D15'; Text = '16.52'; Numeric = $true },
    @{ Cell = 'E15'; Text = '  +1.57%  '; Numeric = $false },
    @{ Cell = 'D16'; Text = '2.730.42'; Numeric = $false },
    @{ Cell = 'E16'; Text = '  +1.14%  '; Numeric = $false },
    @{ Cell = 'D17'; Text = '2.366.41'; Numeric = $false },
    @{ Cell = 'E17'; Text = '  +1.02%  '; Numeric = $false },
    @{ Cell = 'D18'; Text = '7.97'; Numeric = $true },
    @{ Cell = 'E18'; Text = '  +9.67%  '; Numeric = $false },
    @{ Cell = 'D19'; Text = '42.740.19'; Numeric = $false },
    @{ Cell = 'E19'; Text = '  -1.00%  '; Numeric = $false },
    @{ Cell = 'E20'; Text = '  -1.50%  '; Numeric = $false },
    @{ Cell = 'E21'; Text = '  +9.83%  '; Numeric = $false },
    @{ Cell = 'D22'; Text = '76.37'; Numeric = $true },
    @{ Cell = 'E22'; Text = '  +1.36%  '; Numeric = $false },
    @{ Cell = 'D23'; Text = '269.96'; Numeric = $true },
    @{ Cell = 'E23'; Text = '  +6.56%  '; Numeric = $false },
    @{ Cell = 'E24'; Text = '  -10.96%  '; Numeric = $false },
    @{ Cell = 'D25'; Text = '10.12'; Numeric = $true },
    @{ Cell = 'E25'; Text = '  +10.90%  '; Numeric = $false },
    @{ Cell = 'E26'; Text = '  +0.07%  '; Numeric = $false },
    @{ Cell = 'E27'; Text = '  -4.41%  '; Numeric = $false },
    @{ Cell = 'D28'; Text = '23.24'; Numeric = $true },
    @{ Cell = 'E28'; Text = '  +3.77%  '; Numeric = $false },
    @{ Cell = 'E29'; Text = '  -2.60%  '; Numeric = $false },
    @{ Cell = 'D30'; Text = '176.53'; Numeric = $true },
    @{ Cell = 'E30'; Text = '  +1.06%  '; Numeric = $false },
    @{ Cell = 'D31'; Text = '3.09'; Numeric = $true },
    @{ Cell = 'E31'; Text = '  -2.53%  '; Numeric = $false },
    @{ Cell = 'D32'; Text = '0.0904'; Numeric = $true },
    @{ Cell = 'E32'; Text = '  -2.70%  '; Numeric = $false },
    @{ Cell = 'D33'; Text = '35.45'; Numeric = $true },
    @{ Cell = 'E33'; Text = '  -9.95%  '; Numeric = $false },
    @{ Cell = 'D34'; Text = '6.17'; Numeric = $true },
    @{ Cell = 'E34'; Text = '  +2.95%  '; Numeric = $false },
    @{ Cell = 'E35'; Text = '  +0.62%  '; Numeric = $false },
    @{ Cell = 'D36'; Text = '4.63'; Numeric = $true },
    @{ Cell = 'E36'; Text = '  -6.92%  '; Numeric = $false },
    @{ Cell = 'D37'; Text = '2.98'; Numeric = $true },
    @{ Cell = 'E37'; Text = '  +10.81%  '; Numeric = $false },
    @{ Cell = 'E38'; Text = '  -4.44%  '; Numeric = $false },
    @{ Cell = 'D39'; Text = '3.83'; Numeric = $true },
    @{ Cell = 'E39'; Text = '  -7.98%  '; Numeric = $false },
    @{ Cell = 'E40'; Text = '  +1.67%  '; Numeric = $false },
    @{ Cell = 'E41'; Text = '  +4.49%  '; Numeric = $false },
    @{ Cell = 'D42'; Text = '0.236'; Numeric = $true },
    @{ Cell = 'E42'; Text = '  +1.28%  '; Numeric = $false },
    @{ Cell = 'D43'; Text = '70.36'; Numeric = $true },
    @{ Cell = 'E43'; Text = '  -3.66%  '; Numeric = $false },
    @{ Cell = 'E44'; Text = '  -0.08%  '; Numeric = $false },
    @{ Cell = 'D45'; Text = '118.41'; Numeric = $true },
    @{ Cell = 'E45'; Text = '  +6.73%  '; Numeric = $false },
    @{ Cell = 'D46'; Text = '91.62'; Numeric = $true },
    @{ Cell = 'E46'; Text = '  +30.40%  '; Numeric = $false },
    @{ Cell = 'D47'; Text = '11.95'; Numeric = $true },
    @{ Cell = 'E47'; Text = '  -7.34%  '; Numeric = $false },
    @{ Cell = 'D48'; Text = '5.51'; Numeric = $true },
    @{ Cell = 'E48'; Text = '  -2.35%  '; Numeric = $false },
    @{ Cell = 'D49'; Text = '9.22'; Numeric = $true },
    @{ Cell = 'E49'; Text = '  -0.91%  '; Numeric = $false },
    @{ Cell = 'E50'; Text = '  -2.44%  '; Numeric = $false },
    @{ Cell = 'D51'; Text = '1.572.43'; Numeric = $false },
    @{ Cell = 'E51'; Text = '  +5.23%  '; Numeric = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)

    if ($u.Numeric) {
        # Some "Price" cells (column D) look like plain numbers (e.g. "331.98").
        # Assigning that text directly would make Excel coerce it into a Double.
        # Prefix with an apostrophe to force text storage -- exactly like typing
        # it into the cell -- then restore the default "Normal" style so no
        # stray number formatting is left behind on the cell.
        $cell.Value = "'" + $u.Text
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Text
    }
}
